$wb = $excel.ActiveWorkbook

# --- Sheet "Test": the cursor was just left on a different cell, no data changed ---
$ws1 = $wb.Worksheets.Item("Test")
$ws1.Activate() | Out-Null
$ws1.Range("C4").Select() | Out-Null

# --- Sheet "Production": three accounts were migrated over to the new computer ---
$ws2 = $wb.Worksheets.Item("Production")
$ws2.Activate() | Out-Null

# Row 3 - water_earth_402
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "water_earth_402@yahoo.com"
$ws2.Range("C3").Value = "helloGoe234"
$ws2.Range("D3").Value = "helloGoe234"

# Row 4 - daft_williams_405
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "daft_williams_405@yahoo.com"
$ws2.Range("C4").Value = "sdFgsdfg892m45"
$ws2.Range("D4").Value = "sdFgsdfg892m45"

# Row 5 - faiyamR003 (re-added here)
$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "faiyamR003@gmail.com"
$ws2.Range("C5").Value = "beatthestreak3"
$ws2.Range("D5").Value = "beatthestreak3"

# Hyperlink the new e-mail addresses, same as the existing B2 mailto link
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:water_earth_402@yahoo.com")
$ws2.Hyperlinks.Add($ws2.Range("B4"), "mailto:daft_williams_405@yahoo.com")
$ws2.Hyperlinks.Add($ws2.Range("B5"), "mailto:faiyamR003@gmail.com")

# Match the existing "Hyperlink" look used on B2
$ws2.Range("B2").Copy()
$ws2.Range("B3:B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column B needs to be a lot wider now that it holds full e-mail addresses
$ws2.Columns.Item(2).ColumnWidth = 25.3

# Leave the active cell on B4, where the last edit happened
$ws2.Range("B4").Select() | Out-Null
